# Update "想去人数" (want-to-go count) values on the 展览 and 全部类型 sheets
# to reflect newly generated output (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 12803
$ws1.Range("F10").Value = 12719
$ws1.Range("F12").Value = 19
$ws1.Range("F13").Value = 8642
$ws1.Range("F14").Value = 7630
$ws1.Range("F18").Value = 120
$ws1.Range("F22").Value = 373
$ws1.Range("F24").Value = 11

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 12803
$ws4.Range("F11").Value = 12719
$ws4.Range("F13").Value = 19
$ws4.Range("F14").Value = 8642
$ws4.Range("F15").Value = 7630
$ws4.Range("F19").Value = 120
$ws4.Range("F24").Value = 373
$ws4.Range("F26").Value = 11
